# Fix: ipk condition — populate the sample grade rows (row 2 & 3) that were
# previously just empty, specially-styled placeholder cells, and add a new
# third data row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old placeholder cells (C2,E2,F2,C3,E3,F3) carried a one-off "green
# Consolas" style that isn't used anywhere else in the workbook. Drop it so
# the new data uses the default style, same as every other cell.
$ws.Range("C2:C3").ClearFormats()
$ws.Range("E2:F3").ClearFormats()

# A leading apostrophe forces Excel to store a numeric-looking string
# ("2.4", "3.1", ...) as text instead of silently re-parsing it as a
# number. The following ClearFormats() on that same cell drops the
# resulting quote-prefix formatting, leaving the cell on the plain
# default style — matching a shared-string cell with no explicit style.

# Row 2 — nim, nama, nama_sma, penghasilan_orang_tua, prodi, jalur, ip_semester_1..4
$ws.Range("A2").Value = 12345
$ws.Range("B2").Value = "Hesoyam"
$ws.Range("C2").Value = "SMK Negeri 1 Siatasbarita"
$ws.Range("D2").Value = 5000000
$ws.Range("E2").Value = "Teknologi Komputer"
$ws.Range("F2").Value = "UTBK"
$ws.Range("G2").Value = "'2.4"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = "'3.1"
$ws.Range("H2").ClearFormats()
$ws.Range("I2").Value = "'3.2"
$ws.Range("I2").ClearFormats()
$ws.Range("J2").Value = "'3.5"
$ws.Range("J2").ClearFormats()

# Row 3
$ws.Range("A3").Value = 23456
$ws.Range("B3").Value = "Aezakmi"
$ws.Range("C3").Value = "SMA Negri 1 Sidamanik"
$ws.Range("D3").Value = 12000000
$ws.Range("E3").Value = "Manajemen Rekayasa"
$ws.Range("F3").Value = "PMDK"
$ws.Range("G3").Value = "'3.3"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = "'3.2"
$ws.Range("H3").ClearFormats()
$ws.Range("I3").Value = "'3.4"
$ws.Range("I3").ClearFormats()
$ws.Range("J3").Value = "'3.5"
$ws.Range("J3").ClearFormats()

# Row 4 — new row
$ws.Range("A4").Value = 34567
$ws.Range("B4").Value = "Uzumymw"
$ws.Range("C4").Value = "SMAN 4 BINJAI"
$ws.Range("D4").Value = 2000000
$ws.Range("E4").Value = "Teknik Bioproses"
$ws.Range("F4").Value = "USM3"
$ws.Range("G4").Value = "'1.3"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = "'1.2"
$ws.Range("H4").ClearFormats()
$ws.Range("I4").Value = "'2.0"
$ws.Range("I4").ClearFormats()
$ws.Range("J4").Value = "'2.1"
$ws.Range("J4").ClearFormats()

# Widen the columns to fit the new, longer content.
$ws.Columns.Item(1).ColumnWidth = 5.15625
$ws.Columns.Item(2).ColumnWidth = 9.21875
$ws.Columns.Item(3).ColumnWidth = 22.96875
$ws.Columns.Item(5).ColumnWidth = 19.6875
$ws.Columns.Item(6).ColumnWidth = 5.3125
$ws.Columns.Item(11).ColumnWidth = 7.34375

# Final selection left on G2, matching the saved view state.
$ws.Range("G2").Select()
